$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 8 blank rows before row 19 (pushes the old "Total" row from 22 down to 30) ---
$null = $ws.Range("19:26").EntireRow.Insert()
# Remove any formatting that bled into the freshly inserted rows so they start out completely blank
$null = $ws.Range("19:26").EntireRow.Clear()

# --- Fill in the new purchase rows (order chosen so shared-strings are created in the expected sequence) ---
$ws.Range("B19").Value = "Thesis v3.00BOB"
$ws.Range("B20").Value = "Thesis v3.00"
$ws.Range("B21").Value = "Pressure v2.00"
$ws.Range("B22").Value = "Components for Thesis v3.00, Thesis v3.00BOB, and Pressure v2.00"
$ws.Range("C19").Value = "ZuAWe6Su"
$ws.Range("C20").Value = "FFjMkWOU"
$ws.Range("C21").Value = "aF0dGEKE"

$ws.Range("A19").Value = "Osh Park"
$ws.Range("A20").Value = "Osh Park"
$ws.Range("A21").Value = "Osh Park"
$ws.Range("A22").Value = "Sparkfun"
$ws.Range("A23").Value = "Digikey"
$ws.Range("B23").Value = "Components for Thesis v3.00, Thesis v3.00BOB, and Pressure v2.00"
$ws.Range("C22").Value = 829749
$ws.Range("D22").Value = 147.21
$ws.Range("C23").Value = 43539112
$ws.Range("D23").Value = 177.91

# --- Fix up the Total row's formula now that it lives at row 30 ---
$ws.Range("D30").Formula = "=SUM(D2:D27)"

# --- Widen column B ---
$ws.Columns.Item(2).ColumnWidth = 31.86

# --- Update the frozen-pane scroll position / selection ---
$null = $ws.Range("C6").Select()

Write-Host "done"
